$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411

$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931

$ws.Range("A8:XFD9").Select()
